# Update the "Functions" worksheet: add a new "Status" column (F) and
# mark the "Location" class's methods (rows 2-15) as "Done".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Fill "Done" status for the "Location" class rows (rows 2-15) first so the
# "Done" shared string is registered before "Status".
$ws.Range("F2:F15").Value = "Done"

# Add new "Status" column header in F1, matching the bold header style
# already used by A1:E1.
$ws.Range("F1").Value = "Status"
$ws.Range("F1").Font.Bold = $true
